# Update notes on the "laser" lesson.
# 1. Mark "Back" as a grammar-flagged sentence fragment (proofErr gramStart/gramEnd).
# 2. Add a new "Descriptors" row to the programming-terms table.
# 3. Merge the split " - Mission " / "Hard" runs in the 2nd page heading.
# 4. Merge the split "Write a program..." / "You can use " runs in the 2nd page body.

$d = $word.ActiveDocument

# --- Namespace / pkg wrapper helpers ------------------------------------
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# -------------------------------------------------------------------------
# 1) Wrap "Back" in the Directions row with proofErr gramStart/gramEnd.
# -------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Up, Down, Left, Right, Forward, Back")
if ($found) {
    $para = $rng.Paragraphs(1)
    $prng = $para.Range
    $body = '<w:p><w:pPr><w:ind w:firstLine="720"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve">Up, Down, Left, Right, Forward, </w:t></w:r>' `
        + '<w:proofErr w:type="gramStart"/>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>Back</w:t></w:r>' `
        + '<w:proofErr w:type="gramEnd"/>' `
        + '</w:p>'
    $prng.InsertXML($pkgOpen + $body + $pkgClose)
}

# -------------------------------------------------------------------------
# 2) Add the "Descriptors" / "High, Low, Over, Under" row to the table.
# -------------------------------------------------------------------------
$table = $d.Tables(1)
$newRow = $table.Rows.Add()
$newRow.Cells(1).Range.Text = "Descriptors"
$newRow.Cells(2).Range.Text = "High, Low, Over, Under"

# -------------------------------------------------------------------------
# 3) Merge " - Mission " + "Hard" into one run on the 2nd page heading.
# -------------------------------------------------------------------------
$dash = [char]0x2013
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Mission Hard")
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $prng2 = $para2.Range
    $body2 = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:lastRenderedPageBreak/><w:t>Laser Break In</w:t></w:r>' `
        + ('<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> ' + $dash + ' Mission Hard</w:t></w:r>') `
        + '</w:p>'
    $prng2.InsertXML($pkgOpen + $body2 + $pkgClose)
}

# -------------------------------------------------------------------------
# 4) Merge "Write a program...side. " + "You can use " into one run on the
#    2nd page body, leaving the trailing "any words..." run untouched.
# -------------------------------------------------------------------------
$apos = [char]0x2019
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("You can use")
if ($found3) {
    $para3 = $rng3.Paragraphs(1)
    $prng3 = $para3.Range
    $body3 = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve">Write a program to navigate your fellow adventurer through a laser field to get to the prize at the other side. You can use </w:t></w:r>' `
        + ('<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>any words or numbers you' + $apos + 'd like to help you write your program.</w:t></w:r>') `
        + '</w:p>'
    $prng3.InsertXML($pkgOpen + $body3 + $pkgClose)
}

Write-Output "edit.ps1 completed"
